$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Change 1: "sigmoid" paragraph wording tweak
#   "...what a sigmoid function really is, it's that it predict probability..."
#   -> "...what a sigmoid function really is, it predicts probability..."
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "it’s that it predict probability",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "it predicts probability",
    2
) | Out-Null

# ------------------------------------------------------------------
# Change 2: add a new "Feature-Scaling" bullet right after the
# "Imputing" bullet, and move the _GoBack bookmark into the new
# paragraph (right before the closing parenthesis).
# ------------------------------------------------------------------

# Locate the "Imputing" paragraph (the last paragraph of the body).
$imputing = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "*Imputing:*") {
        $imputing = $cand
    }
}

# Drop the old bookmark - it gets re-added further down in the new text.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$fullText = "Q Feature-Scaling: Feature scaling is done in the stage of data pre-processing where you have a data-set with features with gap or differences between them for example you have dataset contain feature in cm and another feature is in meter, this will cause problems with learning algorithms like Logistic regression because the mean of the features in cm might be drastically smaller than the one in meters, so it will cause delay in convergence of gradient descent to reach the minimum point(affecting speed of the learning algorithm)."

# Append the new sentence right onto the end of the "Imputing" paragraph
# text (this keeps it sharing the same, already non-bold, run
# formatting) then split it off into its own paragraph - this avoids
# inheriting stray bold from the document's trailing paragraph mark.
$joinPos = $imputing.Range.End - 1
$d.Range($joinPos, $joinPos).InsertAfter($fullText) | Out-Null
$d.Range($joinPos, $joinPos).InsertParagraphAfter() | Out-Null

$newPara = $d.Paragraphs.Item($imputing.Index + 1)
$base = $newPara.Range.Start

# Bold just the "Feature-Scaling: " label.
$boldStart = $base + 2
$boldEnd = $boldStart + "Feature-Scaling: ".Length
$d.Range($boldStart, $boldEnd).Font.Bold = $true

# Re-insert the _GoBack bookmark right before the final "point(...)"
# parenthesis, matching its original position relative to the text.
$bmOffset = $fullText.IndexOf("eed of the learning algorithm") + "eed of the learning algorithm".Length
$bmPos = $base + $bmOffset
$d.Bookmarks.Add("_GoBack", $d.Range($bmPos, $bmPos)) | Out-Null
